$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text (trim the stray spaces the original file had).
$ws.Range("A1").Value = "BlackBoard"
$ws.Range("B1").Value = "CasePresentation"
$ws.Range("C1").Value = "PPT"

# Widen the columns to comfortably fit the (now un-padded) header text.
$ws.Columns.Item(1).ColumnWidth = 18.02
$ws.Columns.Item(2).ColumnWidth = 21.88
$ws.Columns.Item(3).ColumnWidth = 20.31

# The header row no longer needs to wrap onto multiple lines at the new
# column widths, so it shrinks back down to the sheet's normal row height.
$ws.Rows.Item(1).RowHeight = 19.5

# Move the active selection to C1.
[void]$ws.Range("C1").Select()
